$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The row that held "1304060 - Maria das Graças de Almeida Felipe" underneath
# "Docentes responsáveis:" (row 13, with no A-column label) is removed
# entirely; everything below shifts up by one row.
$ws.Rows(13).Delete()

# After the shift, several B/C cells receive new content so the sheet now
# reads (row: B/C value):
#   10: "1304060 - Maria das Graças de Almeida Felipe"
#   13: "Semestral"
#   15: "01/01/2022"
#   18: "1304060 - Maria das Graças de Almeida Felipe"
#   19: "Duas provas escritas (P1 e P2) distribuídas no semestre."
#   20: "MF=Média finalMF = (P1 + P2) / 2"
#   21: "NF = (MF + PR)/2, onde PR é uma prova de recuperação. Prova de
#        recuperação (PR) para alunos com Média Final maior ou igual a 3,0 e
#        menor do que 5,0. Será considerado aprovado o aluno que tenha obtido
#        Nota Final igual ou maior do que 5,0."
$ws.Range("B10").Value = "1304060 - Maria das Graças de Almeida Felipe"
$ws.Range("C10").Value = "1304060 - Maria das Graças de Almeida Felipe"

$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

$ws.Range("B15").Value = "01/01/2022"
$ws.Range("C15").Value = "01/01/2022"

$ws.Range("B18").Value = "1304060 - Maria das Graças de Almeida Felipe"
$ws.Range("C18").Value = "1304060 - Maria das Graças de Almeida Felipe"

$ws.Range("B19").Value = "Duas provas escritas (P1 e P2) distribuídas no semestre."
$ws.Range("C19").Value = "Duas provas escritas (P1 e P2) distribuídas no semestre."

$ws.Range("B20").Value = "MF=Média finalMF = (P1 + P2) / 2"
$ws.Range("C20").Value = "MF=Média finalMF = (P1 + P2) / 2"

$ws.Range("B21").Value = "NF = (MF + PR)/2, onde PR é uma prova de recuperação. Prova de recuperação (PR) para alunos com Média Final maior ou igual a 3,0 e menor do que 5,0. Será considerado aprovado o aluno que tenha obtido Nota Final igual ou maior do que 5,0."
$ws.Range("C21").Value = "NF = (MF + PR)/2, onde PR é uma prova de recuperação. Prova de recuperação (PR) para alunos com Média Final maior ou igual a 3,0 e menor do que 5,0. Será considerado aprovado o aluno que tenha obtido Nota Final igual ou maior do que 5,0."
